$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the source diff. Cells whose new value is a plain
# decimal number (e.g. "225.73") must be forced to remain stored as text,
# matching the original inline-string cell type, instead of being
# auto-converted to a numeric value by Excel.

$ws.Range("D2").Value = '34.611.37'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.808.96'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.599'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '37.42'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.07%  '
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0683'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '2.073.18'
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").Value = '1.832.49'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '34.605.82'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '0.0₃0776'
$ws.Range("E20").Value = '  -2.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.88%  '
$ws.Range("E28").Value = '  +2.17%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.91%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = '1.366.76'
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.657'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.81%  '
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("E38").Value = '  -5.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0188'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.97%  '
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("E44").Value = '  +4.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("D47").Value = '1.973.46'
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D51").Value = '0.0₆0123'
$ws.Range("E51").Value = '  -5.98%  '
